# Auto-generated edit script: updates cryptos price/volume table
# to match the target commit (GitHub Actions price refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.102.94"
$ws.Range("E2").Value = "'  -0.73%  "
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3").Value = "'1.644.69"
$ws.Range("E3").Value = "'  -0.95%  "
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'217.81"
$ws.Range("E5").Value = "'  -0.81%  "
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("E6").Value = "'  +1.35%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.16%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.21%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.0629"
$ws.Range("E9").Value = "'  +0.44%  "
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'20.00"
$ws.Range("E10").Value = "'  +0.96%  "
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.39%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.875.72"
$ws.Range("E12").Value = "'  -0.92%  "
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("D13").Value = "'1.640.59"
$ws.Range("E13").Value = "'  -1.07%  "
$ws.Range("D13:E13").Style = "Normal"

$ws.Range("D14").Value = "'4.11"
$ws.Range("E14").Value = "'  -2.02%  "
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15").Value = "'0.538"
$ws.Range("E15").Value = "'  +0.72%  "
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'67.35"
$ws.Range("E16").Value = "'  +0.71%  "
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'27.072.34"
$ws.Range("E17").Value = "'  -0.84%  "
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("E18").Value = "'  +0.44%  "
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("D19").Value = "'218.80"
$ws.Range("E19").Value = "'  -1.59%  "
$ws.Range("D19:E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.01%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.85"
$ws.Range("E21").Value = "'  +1.28%  "
$ws.Range("D21:E21").Style = "Normal"

$ws.Range("E22").Value = "'  -0.18%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.22%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -0.55%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'147.88"
$ws.Range("E25").Value = "'  +0.30%  "
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("E26").Value = "'  -0.19%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'7.39"
$ws.Range("E27").Value = "'  -0.62%  "
$ws.Range("D27:E27").Style = "Normal"

$ws.Range("D28").Value = "'0.118"
$ws.Range("E28").Value = "'  -0.71%  "
$ws.Range("D28:E28").Style = "Normal"

$ws.Range("D29").Value = "'15.76"
$ws.Range("E29").Value = "'  -1.48%  "
$ws.Range("D29:E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "'  -1.50%  "
$ws.Range("D30:E30").Style = "Normal"

$ws.Range("E31").Value = "'  -0.79%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.37"
$ws.Range("E32").Value = "'  -0.76%  "
$ws.Range("D32:E32").Style = "Normal"

$ws.Range("E33").Value = "'  +0.93%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.59"
$ws.Range("E34").Value = "'  +1.53%  "
$ws.Range("D34:E34").Style = "Normal"

$ws.Range("D35").Value = "'1.265.04"
$ws.Range("E35").Value = "'  +0.44%  "
$ws.Range("D35:E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.26%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.0177"
$ws.Range("E37").Value = "'  +0.12%  "
$ws.Range("D37:E37").Style = "Normal"

$ws.Range("D38").Value = "'0.541"
$ws.Range("E38").Value = "'  +0.94%  "
$ws.Range("D38:E38").Style = "Normal"

$ws.Range("E39").Value = "'  +1.63%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -0.07%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = "'  -0.52%  "
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = "'  +3.89%  "
$ws.Range("D42:E42").Style = "Normal"

$ws.Range("E43").Value = "'  -0.26%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.786.56"
$ws.Range("E44").Value = "'  -1.11%  "
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("D45").Value = "'62.24"
$ws.Range("E45").Value = "'  +0.85%  "
$ws.Range("D45:E45").Style = "Normal"

$ws.Range("D46").Value = "'91.92"
$ws.Range("E46").Value = "'  -0.76%  "
$ws.Range("D46:E46").Style = "Normal"

$ws.Range("E47").Value = "'  -0.50%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "'  -0.87%  "
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₇0999"
$ws.Range("E49").Value = "'  +8.17%  "
$ws.Range("D49:E49").Style = "Normal"

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.64"
$ws.Range("E50").Value = "'  -0.72%  "
$ws.Range("D50:E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0971"
$ws.Range("E51").Value = "'  -1.11%  "
$ws.Range("D51:E51").Style = "Normal"
